# Deploying to gh-pages from @ codeforIATI/codelists@014ad6fcea0684c5d63b9e87249db0675c32a414
#
# The sector-group codelist columns D/E/F/G were re-ordered from
#   D=codeforiati:group-name, E=codeforiati:category-name,
#   F=codeforiati:group-code,  G=codeforiati:category-code
# to
#   D=codeforiati:category-code, E=codeforiati:group-code,
#   F=codeforiati:group-name,    G=codeforiati:category-name
#
# i.e. for every row, column D swaps with column G, and column E swaps
# with column F. This holds for the header row too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Range("A1").CurrentRegion.Rows.Count
if ($lastRow -lt 235) { $lastRow = 235 }

for ($r = 1; $r -le $lastRow; $r++) {
    $d = $ws.Cells.Item($r, 4).Text
    $e = $ws.Cells.Item($r, 5).Text
    $f = $ws.Cells.Item($r, 6).Text
    $g = $ws.Cells.Item($r, 7).Text

    # New D = old G
    if ($g -match '^[0-9]+$') { $ws.Cells.Item($r, 4).Value = "'" + $g } else { $ws.Cells.Item($r, 4).Value = $g }
    # New E = old F
    if ($f -match '^[0-9]+$') { $ws.Cells.Item($r, 5).Value = "'" + $f } else { $ws.Cells.Item($r, 5).Value = $f }
    # New F = old D
    if ($d -match '^[0-9]+$') { $ws.Cells.Item($r, 6).Value = "'" + $d } else { $ws.Cells.Item($r, 6).Value = $d }
    # New G = old E
    if ($e -match '^[0-9]+$') { $ws.Cells.Item($r, 7).Value = "'" + $e } else { $ws.Cells.Item($r, 7).Value = $e }
}
